$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# 1) Give the AUDREY row (currently row 23) the "closing" bottom-border
#    style that today belongs to the final totals row (row 25) -- once the
#    obsolete rows below are removed, this employee row becomes the very
#    last row of the table and needs that closing border.
# -------------------------------------------------------------------------
$ws.Range("B25:J25").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# -------------------------------------------------------------------------
# 2) Swap the data values (not formats) between row 22 (ALVARO) and row 23
#    (AUDREY) so that AUDREY ends up first (keeping the normal style that
#    was already on row 22) and ALVARO ends up last (keeping the closing
#    style just applied above to row 23).
# -------------------------------------------------------------------------
$alvaro = @($ws.Range("B22").Value2, $ws.Range("C22").Value2, $ws.Range("D22").Value2, $ws.Range("E22").Value2, $ws.Range("F22").Value2, $ws.Range("G22").Value2)
$audrey = @($ws.Range("B23").Value2, $ws.Range("C23").Value2, $ws.Range("D23").Value2, $ws.Range("E23").Value2, $ws.Range("F23").Value2, $ws.Range("G23").Value2)

$ws.Range("B22").Value2 = $audrey[0]
$ws.Range("C22").Value2 = $audrey[1]
$ws.Range("D22").Value2 = $audrey[2]
$ws.Range("E22").Value2 = $audrey[3]
$ws.Range("F22").Value2 = $audrey[4]
$ws.Range("G22").Value2 = $audrey[5]

$ws.Range("B23").Value2 = $alvaro[0]
$ws.Range("C23").Value2 = $alvaro[1]
$ws.Range("D23").Value2 = $alvaro[2]
$ws.Range("E23").Value2 = $alvaro[3]
$ws.Range("F23").Value2 = $alvaro[4]
$ws.Range("G23").Value2 = $alvaro[5]

# -------------------------------------------------------------------------
# 3) Remove the two old NIT subtotal rows (24 and 25) from the bottom up.
# -------------------------------------------------------------------------
$ws.Rows.Item(25).EntireRow.Delete() | Out-Null
$ws.Rows.Item(24).EntireRow.Delete() | Out-Null

# -------------------------------------------------------------------------
# 4) Remove DARLYN ISABEL BARRIOS JIMENEZ's six old late-payment-period
#    rows (16-21) entirely; this shifts AUDREY/ALVARO up to rows 16-17 and
#    the signature block up to rows 22-23.
# -------------------------------------------------------------------------
$ws.Rows("16:21").EntireRow.Delete() | Out-Null

# -------------------------------------------------------------------------
# 5) Refresh the summary figures at the top of the statement.
# -------------------------------------------------------------------------
$ws.Range("E11").Value2 = 75112
$ws.Range("C13").Value2 = 2
$ws.Range("F13").Value2 = 2
